$d = $word.ActiveDocument

# 1) "Porque a marca leva esse nome?" heading: "Porque" -> "Por que"
$r1 = $d.Content
$r1.Find.ClearFormatting()
$r1.Find.Execute("Porque a marca leva esse nome?", $true, $false, $false, $false, $false, $true, 1, $false, "Por que a marca leva esse nome?", 2)

# 2) "Baixo Local de Abrangência" -> "Alto Local de Abrangência"
$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Execute("Baixo Local de Abrangência", $true, $false, $false, $false, $false, $true, 1, $false, "Alto Local de Abrangência", 2)
